$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row Right count (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update the "Total" row Right count (B12): 72 -> 120
$ws.Range("B12").Value = 120

# Update the "Total" row Max text (E12): "71/84" -> "120/140"
$ws.Range("E12").Value = "120/140"
